$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B3").Value = 0.01
$ws1.Range("B6").Value = -263724.5473518896
$ws1.Range("B7").Value = 10477853.63860845
$ws1.Range("B8").Value = 27790152.75161
$ws1.Range("B10").Value = 1684280.80911358

# --- DG Dispatch sheet: populate cells previously 0 with moved Unmet Demand values ---
$ws7 = $wb.Worksheets.Item("DG Dispatch")
$ws7.Range("K2").Value = 220.0898510449805
$ws7.Range("L2").Value = 235.7664149699872
$ws7.Range("N2").Value = 229.4130635965909
$ws7.Range("P2").Value = 231.2329957552695
$ws7.Range("K3").Value = 137.841438974359
$ws7.Range("N3").Value = 131.3417120833333
$ws7.Range("O3").Value = 142.5962444444444
$ws7.Range("L4").Value = 134.8846762812383
$ws7.Range("M4").Value = 138.9257839476051
$ws7.Range("K5").Value = 220.0898510449805
$ws7.Range("N5").Value = 229.4130635965909
$ws7.Range("O5").Value = 230.0982114216867
$ws7.Range("P5").Value = 231.2329957552695
$ws7.Range("K6").Value = 137.841438974359
$ws7.Range("M6").Value = 142.1340339220183
$ws7.Range("P6").Value = 133.9744074143302
$ws7.Range("L7").Value = 134.8846762812383
$ws7.Range("M7").Value = 138.9257839476051
$ws7.Range("N7").Value = 127.6855444652332
$ws7.Range("O7").Value = 138.4565384518428
$ws7.Range("K8").Value = 220.0898510449805
$ws7.Range("L8").Value = 235.7664149699872
$ws7.Range("M8").Value = 230.3462332272727
$ws7.Range("O8").Value = 230.0982114216867
$ws7.Range("K9").Value = 137.841438974359
$ws7.Range("M9").Value = 142.1340339220183
$ws7.Range("O9").Value = 142.5962444444444
$ws7.Range("Q9").Value = 139.9817740860215
$ws7.Range("L10").Value = 134.8846762812383
$ws7.Range("N10").Value = 127.6855444652332
$ws7.Range("M11").Value = 230.3462332272727
$ws7.Range("N11").Value = 229.4130635965909
$ws7.Range("P11").Value = 231.2329957552695
$ws7.Range("K12").Value = 137.841438974359
$ws7.Range("L12").Value = 138.5543797798742
$ws7.Range("M12").Value = 142.1340339220183
$ws7.Range("N12").Value = 131.3417120833333
$ws7.Range("O12").Value = 142.5962444444444
$ws7.Range("P12").Value = 133.9744074143302
$ws7.Range("N13").Value = 127.6855444652332
$ws7.Range("O13").Value = 138.4565384518428
$ws7.Range("M14").Value = 230.3462332272727
$ws7.Range("N14").Value = 229.4130635965909
$ws7.Range("O14").Value = 230.0982114216867
$ws7.Range("P14").Value = 231.2329957552695
$ws7.Range("K15").Value = 137.841438974359
$ws7.Range("M15").Value = 142.1340339220183
$ws7.Range("N15").Value = 131.3417120833333
$ws7.Range("P15").Value = 133.9744074143302
$ws7.Range("Q15").Value = 139.9817740860215
$ws7.Range("M16").Value = 138.9257839476051
$ws7.Range("K17").Value = 220.0898510449805
$ws7.Range("L17").Value = 235.7664149699872
$ws7.Range("N17").Value = 229.4130635965909
$ws7.Range("P17").Value = 231.2329957552695
$ws7.Range("M18").Value = 142.1340339220183
$ws7.Range("P18").Value = 133.9744074143302
$ws7.Range("M20").Value = 230.3462332272727
$ws7.Range("L21").Value = 138.5543797798742
$ws7.Range("O21").Value = 142.5962444444444
$ws7.Range("P21").Value = 133.9744074143302
$ws7.Range("N22").Value = 127.6855444652332
$ws7.Range("L23").Value = 235.7664149699872
$ws7.Range("M23").Value = 230.3462332272727
$ws7.Range("P24").Value = 133.9744074143302
$ws7.Range("N25").Value = 127.6855444652332
$ws7.Range("O25").Value = 138.4565384518428
$ws7.Range("L26").Value = 235.7664149699872
$ws7.Range("M26").Value = 230.3462332272727
$ws7.Range("N26").Value = 229.4130635965909
$ws7.Range("N27").Value = 131.3417120833333
$ws7.Range("Q27").Value = 139.9817740860215
$ws7.Range("L28").Value = 134.8846762812383
$ws7.Range("K29").Value = 220.0898510449805
$ws7.Range("L29").Value = 235.7664149699872
$ws7.Range("M29").Value = 230.3462332272727
$ws7.Range("N29").Value = 229.4130635965909
$ws7.Range("P29").Value = 231.2329957552695
$ws7.Range("K30").Value = 137.841438974359
$ws7.Range("L30").Value = 138.5543797798742
$ws7.Range("M30").Value = 142.1340339220183
$ws7.Range("P30").Value = 133.9744074143302
$ws7.Range("L31").Value = 134.8846762812383
$ws7.Range("M31").Value = 138.9257839476051
$ws7.Range("N31").Value = 127.6855444652332
$ws7.Range("N32").Value = 229.4130635965909
$ws7.Range("K33").Value = 137.841438974359
$ws7.Range("L33").Value = 138.5543797798742
$ws7.Range("M33").Value = 142.1340339220183
$ws7.Range("P33").Value = 133.9744074143302
$ws7.Range("L34").Value = 134.8846762812383
$ws7.Range("M34").Value = 138.9257839476051
$ws7.Range("N34").Value = 127.6855444652332
$ws7.Range("O34").Value = 138.4565384518428
$ws7.Range("K35").Value = 220.0898510449805
$ws7.Range("N35").Value = 229.4130635965909
$ws7.Range("K36").Value = 137.841438974359
$ws7.Range("L36").Value = 138.5543797798742
$ws7.Range("M36").Value = 142.1340339220183
$ws7.Range("N36").Value = 131.3417120833333
$ws7.Range("P36").Value = 133.9744074143302
$ws7.Range("L37").Value = 134.8846762812383
$ws7.Range("M37").Value = 138.9257839476051
$ws7.Range("O37").Value = 138.4565384518428
$ws7.Range("K38").Value = 220.0898510449805
$ws7.Range("L38").Value = 235.7664149699872
$ws7.Range("N38").Value = 229.4130635965909
$ws7.Range("K39").Value = 137.841438974359
$ws7.Range("Q39").Value = 139.9817740860215
$ws7.Range("L40").Value = 134.8846762812383
$ws7.Range("M40").Value = 138.9257839476051
$ws7.Range("N40").Value = 127.6855444652332
$ws7.Range("M41").Value = 230.3462332272727
$ws7.Range("O41").Value = 230.0982114216867
$ws7.Range("P41").Value = 231.2329957552695
$ws7.Range("K42").Value = 137.841438974359
$ws7.Range("L42").Value = 138.5543797798742
$ws7.Range("Q42").Value = 139.9817740860215
$ws7.Range("L43").Value = 134.8846762812383
$ws7.Range("N43").Value = 127.6855444652332
$ws7.Range("O43").Value = 138.4565384518428
$ws7.Range("N44").Value = 229.4130635965909
$ws7.Range("L45").Value = 138.5543797798742
$ws7.Range("M45").Value = 142.1340339220183
$ws7.Range("Q45").Value = 139.9817740860215
$ws7.Range("M46").Value = 138.9257839476051

# --- Unmet Demand sheet: zero out the cells that moved to DG Dispatch ---
$ws16 = $wb.Worksheets.Item("Unmet Demand")
$ws16.Range("K2").Value = 0
$ws16.Range("L2").Value = 0
$ws16.Range("N2").Value = 0
$ws16.Range("P2").Value = 0
$ws16.Range("K3").Value = 0
$ws16.Range("N3").Value = 0
$ws16.Range("O3").Value = 0
$ws16.Range("L4").Value = 0
$ws16.Range("M4").Value = 0
$ws16.Range("K5").Value = 0
$ws16.Range("N5").Value = 0
$ws16.Range("O5").Value = 0
$ws16.Range("P5").Value = 0
$ws16.Range("K6").Value = 0
$ws16.Range("M6").Value = 0
$ws16.Range("P6").Value = 0
$ws16.Range("L7").Value = 0
$ws16.Range("M7").Value = 0
$ws16.Range("N7").Value = 0
$ws16.Range("O7").Value = 0
$ws16.Range("K8").Value = 0
$ws16.Range("L8").Value = 0
$ws16.Range("M8").Value = 0
$ws16.Range("O8").Value = 0
$ws16.Range("K9").Value = 0
$ws16.Range("M9").Value = 0
$ws16.Range("O9").Value = 0
$ws16.Range("Q9").Value = 0
$ws16.Range("L10").Value = 0
$ws16.Range("N10").Value = 0
$ws16.Range("M11").Value = 0
$ws16.Range("N11").Value = 0
$ws16.Range("P11").Value = 0
$ws16.Range("K12").Value = 0
$ws16.Range("L12").Value = 0
$ws16.Range("M12").Value = 0
$ws16.Range("N12").Value = 0
$ws16.Range("O12").Value = 0
$ws16.Range("P12").Value = 0
$ws16.Range("N13").Value = 0
$ws16.Range("O13").Value = 0
$ws16.Range("M14").Value = 0
$ws16.Range("N14").Value = 0
$ws16.Range("O14").Value = 0
$ws16.Range("P14").Value = 0
$ws16.Range("K15").Value = 0
$ws16.Range("M15").Value = 0
$ws16.Range("N15").Value = 0
$ws16.Range("P15").Value = 0
$ws16.Range("Q15").Value = 0
$ws16.Range("M16").Value = 0
$ws16.Range("K17").Value = 0
$ws16.Range("L17").Value = 0
$ws16.Range("N17").Value = 0
$ws16.Range("P17").Value = 0
$ws16.Range("M18").Value = 0
$ws16.Range("P18").Value = 0
$ws16.Range("M20").Value = 0
$ws16.Range("L21").Value = 0
$ws16.Range("O21").Value = 0
$ws16.Range("P21").Value = 0
$ws16.Range("N22").Value = 0
$ws16.Range("L23").Value = 0
$ws16.Range("M23").Value = 0
$ws16.Range("P24").Value = 0
$ws16.Range("N25").Value = 0
$ws16.Range("O25").Value = 0
$ws16.Range("L26").Value = 0
$ws16.Range("M26").Value = 0
$ws16.Range("N26").Value = 0
$ws16.Range("N27").Value = 0
$ws16.Range("Q27").Value = 0
$ws16.Range("L28").Value = 0
$ws16.Range("K29").Value = 0
$ws16.Range("L29").Value = 0
$ws16.Range("M29").Value = 0
$ws16.Range("N29").Value = 0
$ws16.Range("P29").Value = 0
$ws16.Range("K30").Value = 0
$ws16.Range("L30").Value = 0
$ws16.Range("M30").Value = 0
$ws16.Range("P30").Value = 0
$ws16.Range("L31").Value = 0
$ws16.Range("M31").Value = 0
$ws16.Range("N31").Value = 0
$ws16.Range("N32").Value = 0
$ws16.Range("K33").Value = 0
$ws16.Range("L33").Value = 0
$ws16.Range("M33").Value = 0
$ws16.Range("P33").Value = 0
$ws16.Range("L34").Value = 0
$ws16.Range("M34").Value = 0
$ws16.Range("N34").Value = 0
$ws16.Range("O34").Value = 0
$ws16.Range("K35").Value = 0
$ws16.Range("N35").Value = 0
$ws16.Range("K36").Value = 0
$ws16.Range("L36").Value = 0
$ws16.Range("M36").Value = 0
$ws16.Range("N36").Value = 0
$ws16.Range("P36").Value = 0
$ws16.Range("L37").Value = 0
$ws16.Range("M37").Value = 0
$ws16.Range("O37").Value = 0
$ws16.Range("K38").Value = 0
$ws16.Range("L38").Value = 0
$ws16.Range("N38").Value = 0
$ws16.Range("K39").Value = 0
$ws16.Range("Q39").Value = 0
$ws16.Range("L40").Value = 0
$ws16.Range("M40").Value = 0
$ws16.Range("N40").Value = 0
$ws16.Range("M41").Value = 0
$ws16.Range("O41").Value = 0
$ws16.Range("P41").Value = 0
$ws16.Range("K42").Value = 0
$ws16.Range("L42").Value = 0
$ws16.Range("Q42").Value = 0
$ws16.Range("L43").Value = 0
$ws16.Range("N43").Value = 0
$ws16.Range("O43").Value = 0
$ws16.Range("N44").Value = 0
$ws16.Range("L45").Value = 0
$ws16.Range("M45").Value = 0
$ws16.Range("Q45").Value = 0
$ws16.Range("M46").Value = 0

# --- Household Surplus sheet ---
$ws17 = $wb.Worksheets.Item("Household Surplus")
$ws17.Range("B2").Value = 172147.4564623187
$ws17.Range("B3").Value = 182847.5994019398
$ws17.Range("B4").Value = 182987.3945782901
$ws17.Range("B5").Value = 171165.9949734709
$ws17.Range("B6").Value = 187632.7892677333
$ws17.Range("B7").Value = 150272.5348460527
$ws17.Range("B8").Value = 69227.33336802496
$ws17.Range("B9").Value = 87277.62269806072
$ws17.Range("B10").Value = 124436.7205980396
$ws17.Range("B11").Value = 220723.4332581452
$ws17.Range("B12").Value = 96997.22070799567
$ws17.Range("B13").Value = 132539.0888193136
$ws17.Range("B14").Value = 134834.408395709
$ws17.Range("B15").Value = 146147.4231337193
$ws17.Range("B16").Value = 69991.64848464866

# --- Costs and Revenues sheet ---
$ws2 = $wb.Worksheets.Item("Costs and Revenues")
$ws2.Range("B4").Value = 2424.612062849559
$ws2.Range("C4").Value = 2575.318301435772
$ws2.Range("D4").Value = 2577.287247581552
$ws2.Range("E4").Value = 2410.788661598182
$ws2.Range("F4").Value = 2642.715341799061
$ws2.Range("G4").Value = 2116.514575296517
$ws2.Range("H4").Value = 975.0328643383805
$ws2.Range("I4").Value = 1229.262291521983
$ws2.Range("J4").Value = 1752.629867578022
$ws2.Range("K4").Value = 3108.780750114724
$ws2.Range("L4").Value = 1366.158038140784
$ws2.Range("M4").Value = 1866.747729849489
$ws2.Range("N4").Value = 1899.076174587451
$ws2.Range("O4").Value = 2058.414410334074
$ws2.Range("P4").Value = 985.7978659809678
$ws2.Range("B6").Value = -54153.64424660709
$ws2.Range("C6").Value = -54153.64424660708
$ws2.Range("D6").Value = -54153.64424660709
$ws2.Range("E6").Value = -20526.04424660708
$ws2.Range("F6").Value = -20526.04424660708
$ws2.Range("G6").Value = -20526.04424660709
$ws2.Range("H6").Value = -20526.04424660708
$ws2.Range("I6").Value = -20526.04424660709
$ws2.Range("J6").Value = -20526.04424660708
$ws2.Range("K6").Value = -20526.04424660708
$ws2.Range("L6").Value = -20526.04424660708
$ws2.Range("M6").Value = -20526.04424660708
$ws2.Range("N6").Value = -20526.04424660708
$ws2.Range("O6").Value = -20526.04424660708
$ws2.Range("P6").Value = -20526.04424660708
